$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 360.6
$ws.Range("J5").Value = 549.6667
$ws.Range("L5").Value = 549.6667
$ws.Range("N5").Value = -779.6667
$ws.Range("H55").Value = 227.375
$ws.Range("J55").Value = 91.666664
$ws.Range("L55").Value = 91.666664
$ws.Range("N55").Value = -519.666664
$ws.Range("H76").Value = 83337030
$ws.Range("I76").Value = 111114640
$ws.Range("J76").Value = 4199.6665
$ws.Range("K76").Value = 111114640
$ws.Range("L76").Value = 4199.6665
$ws.Range("M76").Value = -111114325
$ws.Range("N76").Value = -4829.6665
$ws.Range("H79").Value = 83337030
$ws.Range("I79").Value = 111114640
$ws.Range("J79").Value = 4199.6665
$ws.Range("K79").Value = 111114640
$ws.Range("L79").Value = 4199.6665
$ws.Range("M79").Value = -111113548
$ws.Range("N79").Value = -6383.6665
$ws.Range("H125").Value = 4096.923
$ws.Range("I125").Value = 2677.8333
$ws.Range("J125").Value = 5313.2856
$ws.Range("K125").Value = 24100.4997
$ws.Range("L125").Value = 47819.5704
$ws.Range("M125").Value = -21640.4997
$ws.Range("N125").Value = -52739.5704
$ws.Range("H134").Value = 34311.812
$ws.Range("J134").Value = 34311.812
$ws.Range("L134").Value = 34311.812
$ws.Range("N134").Value = -44451.812
$ws.Range("H137").Value = 1179.5385
$ws.Range("I137").Value = 967.3929000000001
$ws.Range("K137").Value = 2902.1787
$ws.Range("M137").Value = -352.1787000000004

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3515.25
$ws.Range("J45").Value = 3742.5557
$ws.Range("L45").Value = 3742.5557
$ws.Range("N45").Value = -4496.5557
$ws.Range("H63").Value = 1291
$ws.Range("I63").Value = 1089.1428
$ws.Range("K63").Value = 1089.1428
$ws.Range("M63").Value = -403.1428000000001
$ws.Range("H66").Value = 1291
$ws.Range("I66").Value = 1089.1428
$ws.Range("K66").Value = 5445.714
$ws.Range("M66").Value = -2013.714
$ws.Range("H97").Value = 3306.68
$ws.Range("I97").Value = 776.4
$ws.Range("K97").Value = 776.4
$ws.Range("M97").Value = -280.4
$ws.Range("H122").Value = 2884.125
$ws.Range("I122").Value = 2846.5
$ws.Range("J122").Value = 2997
$ws.Range("K122").Value = 8539.5
$ws.Range("L122").Value = 8991
$ws.Range("M122").Value = -6089.5
$ws.Range("N122").Value = -13891
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null
$ws.Range("H132").Value = 5143.44
$ws.Range("I132").Value = 4981.273
$ws.Range("K132").Value = 14943.819
$ws.Range("M132").Value = -12413.819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 880.55554
$ws.Range("I99").Value = 847.3333
$ws.Range("K99").Value = 847.3333
$ws.Range("M99").Value = 650.6667
$ws.Range("H105").Value = 3314.5925
$ws.Range("I105").Value = 2228.4707
$ws.Range("K105").Value = 2228.4707
$ws.Range("M105").Value = -481.4706999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2522.4565
$ws.Range("I31").Value = 2040.1154
$ws.Range("J31").Value = 3149.5
$ws.Range("K31").Value = 2040.1154
$ws.Range("L31").Value = 3149.5
$ws.Range("M31").Value = -1745.1154
$ws.Range("N31").Value = -3739.5
$ws.Range("H34").Value = 2522.4565
$ws.Range("I34").Value = 2040.1154
$ws.Range("J34").Value = 3149.5
$ws.Range("K34").Value = 2040.1154
$ws.Range("L34").Value = 3149.5
$ws.Range("M34").Value = -1838.1154
$ws.Range("N34").Value = -3553.5
$ws.Range("H39").Value = 16762.375
$ws.Range("I39").Value = 16762.375
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 16762.375
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -16371.375
$ws.Range("N39").Value = $null
$ws.Range("H49").Value = 16762.375
$ws.Range("I49").Value = 16762.375
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 16762.375
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -16580.375
$ws.Range("N49").Value = $null
$ws.Range("H62").Value = 13553.1
$ws.Range("I62").Value = 3598.4
$ws.Range("K62").Value = 3598.4
$ws.Range("M62").Value = -2974.4
$ws.Range("H65").Value = 13553.1
$ws.Range("I65").Value = 3598.4
$ws.Range("K65").Value = 17992
$ws.Range("M65").Value = -14872
$ws.Range("H105").Value = 1584.9231
$ws.Range("I105").Value = 1645.909
$ws.Range("J105").Value = 1249.5
$ws.Range("K105").Value = 1645.909
$ws.Range("L105").Value = 1249.5
$ws.Range("M105").Value = 101.0909999999999
$ws.Range("N105").Value = -4743.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 774.1667
$ws.Range("I2").Value = 249.66667
$ws.Range("J2").Value = 1298.6666
$ws.Range("K2").Value = 1498.00002
$ws.Range("L2").Value = 7791.9996
$ws.Range("M2").Value = -1385.00002
$ws.Range("N2").Value = -8017.9996
$ws.Range("H16").Value = 4136.6665
$ws.Range("I16").Value = 705.5
$ws.Range("K16").Value = 2116.5
$ws.Range("M16").Value = -1943.5
$ws.Range("H68").Value = 1816.6666
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 3133.3333
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 9399.999899999999
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -11021.9999
$ws.Range("H70").Value = 6282.4736
$ws.Range("J70").Value = 6862
$ws.Range("L70").Value = 20586
$ws.Range("N70").Value = -21216
$ws.Range("H71").Value = 1816.6666
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 3133.3333
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 28199.9997
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -36311.9997
$ws.Range("H73").Value = 6282.4736
$ws.Range("J73").Value = 6862
$ws.Range("L73").Value = 20586
$ws.Range("N73").Value = -22770
$ws.Range("H109").Value = 549.5
$ws.Range("I109").Value = 549.5
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1648.5
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -608.5
$ws.Range("N109").Value = $null
$ws.Range("H113").Value = 863.51514
$ws.Range("I113").Value = 551.25
$ws.Range("J113").Value = 963.4400000000001
$ws.Range("K113").Value = 1653.75
$ws.Range("L113").Value = 2890.32
$ws.Range("M113").Value = 516.25
$ws.Range("N113").Value = -7230.32
$ws.Range("H132").Value = 4415.88
$ws.Range("I132").Value = 2228.6
$ws.Range("J132").Value = 4962.7
$ws.Range("K132").Value = 20057.4
$ws.Range("L132").Value = 44664.3
$ws.Range("M132").Value = -17527.4
$ws.Range("N132").Value = -49724.3

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3833.8696
$ws.Range("I80").Value = 3472.3845
$ws.Range("K80").Value = 3472.3845
$ws.Range("M80").Value = -2474.3845
$ws.Range("H83").Value = 3833.8696
$ws.Range("I83").Value = 3472.3845
$ws.Range("K83").Value = 17361.9225
$ws.Range("M83").Value = -12369.9225
$ws.Range("H102").Value = 11811.28
$ws.Range("I102").Value = 12220.917
$ws.Range("K102").Value = 12220.917
$ws.Range("M102").Value = -10598.917
$ws.Range("H113").Value = 4870.75
$ws.Range("I113").Value = 2999.3333
$ws.Range("K113").Value = 2999.3333
$ws.Range("M113").Value = -829.3332999999998
$ws.Range("H135").Value = 87827.07000000001
$ws.Range("I135").Value = 79999
$ws.Range("J135").Value = 88429.234
$ws.Range("K135").Value = 79999
$ws.Range("L135").Value = 88429.234
$ws.Range("M135").Value = -74929
$ws.Range("N135").Value = -98569.234

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2751.4375
$ws.Range("I40").Value = 2571.182
$ws.Range("K40").Value = 2571.182
$ws.Range("M40").Value = -2435.182
$ws.Range("H46").Value = 1889.4
$ws.Range("I46").Value = 1599.5
$ws.Range("K46").Value = 1599.5
$ws.Range("M46").Value = -1411.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
$ws.Range("H100").Value = 2575.1667
$ws.Range("I100").Value = 1947.1666
$ws.Range("K100").Value = 1947.1666
$ws.Range("M100").Value = -1406.1666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2092.1738
$ws.Range("I126").Value = 1848.5883
$ws.Range("K126").Value = 5545.7649
$ws.Range("M126").Value = -3075.7649
